$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.957.48"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "3.919.05"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "608.29"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "170.14"
$ws.Range("E6").Value = "  +4.84%  "
$ws.Range("D7").Value = "3.918.97"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E13").Value = "  +5.04%  "
$ws.Range("D14").Value = "38.39"
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("D15").Value = "4.577.89"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "3.927.03"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "69.971.15"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "18.73"
$ws.Range("E18").Value = "  +9.68%  "
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "494.16"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "0.747"
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").Value = "85.74"
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").Value = "12.34"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").Value = "4.070.63"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").Value = "2.44"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").Value = "7.85"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").Value = "32.20"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "3.883.05"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "6.14"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "1.04"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").Value = "3.31"
$ws.Range("E40").Value = "  +11.89%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.330"
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "2.13"
$ws.Range("E43").Value = "  +7.67%  "
$ws.Range("D44").Value = "438.86"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "48.24"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").Value = "8.67"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").Value = "40.94"
$ws.Range("E49").Value = "  +5.91%  "
$ws.Range("D50").Value = "0.000274"
$ws.Range("E50").Value = "  +20.54%  "
$ws.Range("D51").Value = "143.41"
$ws.Range("E51").Value = "  -0.05%  "
